$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

# Update "Marking" row Right-answer count
$ws.Range("B11").Value = 5

# Update "Total" row Right-answer count
$ws.Range("B12").Value = 70

# Update Correct/Total marks text
$ws.Range("E12").Value = "70/140"
